# "Generate Report for Handback" - fills in handback columns for each
# localized-file row on the zh-cn and de-de sheets, marks the overview /
# per-language Status as handed-back-in-sync, widens the columns that now
# hold the longer file names, and hyperlinks the newly-populated
# "Latest Target File" cells to the source markdown files (same target as
# column A's existing hyperlink for that row).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Handed back: in sync with en-US"

$url5f = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/752cf5ad112bc7a7458ecb013dc89f9008629b86/e2e/5f1db785-3cd5-4a7d-88ec-441d2527faf8.md"
$urlf1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/752cf5ad112bc7a7458ecb013dc89f9008629b86/e2e/f1fde0df-f0b6-4ce6-9e7c-811e0aa149bf.md"

$md5f  = "5f1db785-3cd5-4a7d-88ec-441d2527faf8.md"
$mdf1  = "f1fde0df-f0b6-4ce6-9e7c-811e0aa149bf.md"

$zhxlf5f = "5f1db785-3cd5-4a7d-88ec-441d2527faf8.467e21cac036376dab117a82566d6c8c38adbc7c.zh-cn.xlf"
$zhxlff1 = "f1fde0df-f0b6-4ce6-9e7c-811e0aa149bf.8abda619694070839db5548b71793569ab5af701.zh-cn.xlf"
$dexlf5f = "5f1db785-3cd5-4a7d-88ec-441d2527faf8.467e21cac036376dab117a82566d6c8c38adbc7c.de-de.xlf"
$dexlff1 = "f1fde0df-f0b6-4ce6-9e7c-811e0aa149bf.8abda619694070839db5548b71793569ab5af701.de-de.xlf"

$zhHandback = "2016-08-21 18:35:46"
$deHandback = "2016-08-21 18:35:52"

# ---- Overview sheet: status text + widen the two language status columns ----
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Columns.Item(5).ColumnWidth = 29.1666667
$ws1.Columns.Item(6).ColumnWidth = 29.1666667

# ---- zh-cn sheet ----
$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws2.Range("I2").Value = $md5f
$ws2.Range("I2").Style = "Hyperlink"
$ws2.Range("J2").Value = $zhxlf5f
$ws2.Range("K2").Value = $zhHandback

$ws2.Range("I3").Value = $mdf1
$ws2.Range("I3").Style = "Hyperlink"
$ws2.Range("J3").Value = $zhxlff1
$ws2.Range("K3").Value = $zhHandback

$ws2.Hyperlinks.Add($ws2.Range("I2"), $url5f, "", "", $md5f) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), $urlf1, "", "", $mdf1) | Out-Null

$ws2.Columns.Item(3).ColumnWidth = 29.1666667
$ws2.Columns.Item(9).ColumnWidth = 39.1666667
$ws2.Columns.Item(10).ColumnWidth = 39.1666667

# ---- de-de sheet ----
$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

$ws3.Range("I2").Value = $md5f
$ws3.Range("I2").Style = "Hyperlink"
$ws3.Range("J2").Value = $dexlf5f
$ws3.Range("K2").Value = $deHandback

$ws3.Range("I3").Value = $mdf1
$ws3.Range("I3").Style = "Hyperlink"
$ws3.Range("J3").Value = $dexlff1
$ws3.Range("K3").Value = $deHandback

$ws3.Hyperlinks.Add($ws3.Range("I2"), $url5f, "", "", $md5f) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), $urlf1, "", "", $mdf1) | Out-Null

$ws3.Columns.Item(3).ColumnWidth = 29.1666667
$ws3.Columns.Item(9).ColumnWidth = 39.1666667
$ws3.Columns.Item(10).ColumnWidth = 39.1666667
